$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.151.25"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.789.09"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'227.08"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'0.547"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "'0.0691"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "2.047.99"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").Value = "1.794.61"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "34.098.67"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "'0.621"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'4.18"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "'68.14"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'245.68"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "0.0₃0779"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'10.88"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "'161.52"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "'16.34"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'0.0517"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "'3.66"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'3.61"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "1.456.81"
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = "  +9.74%  "
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("D39").Value = "'1.04"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "'80.36"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "'0.920"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'13.52"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("D45").Value = "'6.05"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "0.0₆0136"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.949.03"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'106.25"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("E51").Value = "  +0.00%  "
